# "re-doing the group typology using quantile heatmap"
# The "groups" sheet previously held several repeated 3-column blocks
# (G:H, J:K, M:N, P:Q, S:T) of manually colour-coded "heatmap" values.
# This edit clears all of that out, leaving only the first two data
# columns (A/B and D/E) with the re-computed group typology numbers,
# removes the custom fill colours from the whole grid (back to "No Fill"),
# and moves the active selection to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("groups")
$ws.Activate()

# ---- 1. Update the surviving data columns (A, B, D, E) for rows 3-7 ----
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 5
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 7
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 5

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 8
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 6

$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = 9
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 8

# ---- 2. Clear out the old repeated heatmap blocks (C and F:U, rows 3-9) ----
$ws.Range("C3:C9").ClearContents()
$ws.Range("F3:U9").ClearContents()

# ---- 3. Strip the old manual "heatmap" fill colours from the whole grid ----
$rng = $ws.Range("A3:U9")
$rng.ClearFormats()
$rng.Interior.ColorIndex = -4142

# ---- 4. Move the selection, as seen in the saved file ----
$ws.Range("G3").Select()
